{"js": "// Remove the \"rad\" prefix/substring from a handful of radiometric terms\n// (\"radiance\" -> \"iance\", \"irradiance\" -> \"iriance\", \"radiant\" -> \"iant\",\n// \"mrad\" -> \"m\") throughout the document body. This mirrors a blind\n// literal find-and-replace of \"rad\" -> \"\" that the author ran across the\n// file (see commit message: \"Remove 'rad' prefix from newer routines in\n// PsychRadiometric\").\n\nconst body = context.document.body;\n\n// Case-sensitive, not whole-word, so it matches the substring \"rad\"\n// wherever it occurs (radiance, irradiance, radiant, mrad, ...).\nconst results = body.search(\"rad\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"radiance\"\n$find.Replacement.Text = \"iance\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
